$wb = $excel.ActiveWorkbook

# "recommended_rank" (column B, rows 2-8): identical shared-string table entries
# reused across every sheet, so every sheets B2:B8 must be rewritten to release
# the old shared-string slot and land on the same new text/order as the source data.
$recommendedRank = @(
    '(''RandomForest'', (2.151429674410282, 0.7594959733070695))'
    '(''Knn10'', (1.926330030521544, 0.7319053522713357))'
    '(''DecisionTree'', (1.8917505570038493, 0.7540653475337438))'
    '(''Knn5'', (1.6857194765779724, 0.7128458097177036))'
    '(''LDA'', (1.2373644558632935, 0.6753846338630222))'
    '(''Knn1'', (1.0809838643784946, 0.6472733703791613))'
    '(''NaiveBayes'', (0.14278983669165277, 0.5294468484269127))'
)

$sheetNames = @("abalone", "adult", "banknote", "car", "chess1", "chess2", "contraceptive")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $recommendedRank.Length; $i++) {
        $ws.Cells.Item($i + 2, 2).Value = $recommendedRank[$i]
    }
}

# "ideal_rank" (column C, rows 2-8) and "sperman_coef" (column D) are per-sheet.
$ws = $wb.Worksheets.Item("abalone")
$ws.Range("C2").Value = '(''DecisionTree'', (1.2123171124514067, 0.7540653475337438))'
$ws.Range("D2").Value = 0.8571428571428572
$ws.Range("C3").Value = '(''RandomForest'', (1.1719507301248393, 0.7594959733070695))'
$ws.Range("D3").Value = 0.8571428571428572
$ws.Range("C4").Value = '(''Knn10'', (1.1107397792349358, 0.7319053522713357))'
$ws.Range("D4").Value = 0.8571428571428572
$ws.Range("C5").Value = '(''Knn5'', (1.078780913976465, 0.7128458097177036))'
$ws.Range("D5").Value = 0.8571428571428572
$ws.Range("C6").Value = '(''Knn1'', (0.9315950378153693, 0.6472733703791613))'
$ws.Range("D6").Value = 0.8571428571428572
$ws.Range("C7").Value = '(''LDA'', (0.8942340309054552, 0.6753846338630222))'
$ws.Range("D7").Value = 0.8571428571428572
$ws.Range("C8").Value = '(''NaiveBayes'', (0.6387237856126603, 0.5294468484269127))'
$ws.Range("D8").Value = 0.8571428571428572

$ws = $wb.Worksheets.Item("adult")
$ws.Range("C2").Value = '(''DecisionTree'', (1.2113278699772818, 0.7540653475337438))'
$ws.Range("D2").Value = 0.8571428571428572
$ws.Range("C3").Value = '(''RandomForest'', (1.2102631078639372, 0.7594959733070695))'
$ws.Range("D3").Value = 0.8571428571428572
$ws.Range("C4").Value = '(''Knn10'', (1.1672089703031432, 0.7319053522713357))'
$ws.Range("D4").Value = 0.8571428571428572
$ws.Range("C5").Value = '(''Knn5'', (1.1190590580108295, 0.7128458097177036))'
$ws.Range("D5").Value = 0.8571428571428572
$ws.Range("C6").Value = '(''Knn1'', (0.9497811465864677, 0.6472733703791613))'
$ws.Range("D6").Value = 0.8571428571428572
$ws.Range("C7").Value = '(''LDA'', (0.9464924661090054, 0.6753846338630222))'
$ws.Range("D7").Value = 0.8571428571428572
$ws.Range("C8").Value = '(''NaiveBayes'', (0.5600623026359991, 0.5294468484269127))'
$ws.Range("D8").Value = 0.8571428571428572

$ws = $wb.Worksheets.Item("banknote")
$ws.Range("C2").Value = '(''RandomForest'', (1.2170645838379592, 0.7594959733070695))'
$ws.Range("C3").Value = '(''DecisionTree'', (1.2117369461265741, 0.7540653475337438))'
$ws.Range("C4").Value = '(''Knn10'', (1.1603714022565117, 0.7319053522713357))'
$ws.Range("C5").Value = '(''Knn5'', (1.1082003281710633, 0.7128458097177036))'
$ws.Range("C6").Value = '(''LDA'', (0.9534094428581625, 0.6753846338630222))'
$ws.Range("C7").Value = '(''Knn1'', (0.9305341827483555, 0.6472733703791613))'
$ws.Range("C8").Value = '(''NaiveBayes'', (0.5811958730176384, 0.5294468484269127))'

$ws = $wb.Worksheets.Item("car")
$ws.Range("C2").Value = '(''RandomForest'', (1.2048538766631265, 0.7594959733070695))'
$ws.Range("C3").Value = '(''DecisionTree'', (1.1962394900874063, 0.7540653475337438))'
$ws.Range("C4").Value = '(''Knn10'', (1.1571285853165822, 0.7319053522713357))'
$ws.Range("C5").Value = '(''Knn5'', (1.1157091386800295, 0.7128458097177036))'
$ws.Range("C6").Value = '(''Knn1'', (0.9589110498353272, 0.6472733703791613))'
$ws.Range("C7").Value = '(''LDA'', (0.9516994161223318, 0.6753846338630222))'
$ws.Range("C8").Value = '(''NaiveBayes'', (0.5724149429815409, 0.5294468484269127))'

$ws = $wb.Worksheets.Item("chess1")
$ws.Range("C2").Value = '(''RandomForest'', (1.20142055609216, 0.7594959733070695))'
$ws.Range("C3").Value = '(''DecisionTree'', (1.1942612738535803, 0.7540653475337438))'
$ws.Range("C4").Value = '(''Knn10'', (1.1539285184560222, 0.7319053522713357))'
$ws.Range("C5").Value = '(''Knn5'', (1.1003799029232246, 0.7128458097177036))'
$ws.Range("C6").Value = '(''LDA'', (0.9450467796729798, 0.6753846338630222))'
$ws.Range("C7").Value = '(''Knn1'', (0.9358041383162409, 0.6472733703791613))'
$ws.Range("C8").Value = '(''NaiveBayes'', (0.6084430371476631, 0.5294468484269127))'

$ws = $wb.Worksheets.Item("chess2")
$ws.Range("C2").Value = '(''RandomForest'', (0.9714381179467713, 0.7594959733070695))'
$ws.Range("D2").Value = 0.8571428571428572
$ws.Range("C3").Value = '(''Knn10'', (0.9516558424331265, 0.7319053522713357))'
$ws.Range("D3").Value = 0.8571428571428572
$ws.Range("C4").Value = '(''LDA'', (0.947583130314012, 0.6753846338630222))'
$ws.Range("D4").Value = 0.8571428571428572
$ws.Range("C5").Value = '(''Knn5'', (0.9090488404151519, 0.7128458097177036))'
$ws.Range("D5").Value = 0.8571428571428572
$ws.Range("C6").Value = '(''DecisionTree'', (0.9013743553037171, 0.7540653475337438))'
$ws.Range("D6").Value = 0.8571428571428572
$ws.Range("C7").Value = '(''Knn1'', (0.8193326396879507, 0.6472733703791613))'
$ws.Range("D7").Value = 0.8571428571428572
$ws.Range("C8").Value = '(''NaiveBayes'', (0.6823967374042752, 0.5294468484269127))'
$ws.Range("D8").Value = 0.8571428571428572

$ws = $wb.Worksheets.Item("contraceptive")
$ws.Range("C2").Value = '(''DecisionTree'', (1.2202819478702798, 0.7540653475337438))'
$ws.Range("D2").Value = 0.8928571428571429
$ws.Range("C3").Value = '(''RandomForest'', (1.2074433108774394, 0.7594959733070695))'
$ws.Range("D3").Value = 0.8928571428571429
$ws.Range("C4").Value = '(''Knn10'', (1.1524296008964432, 0.7319053522713357))'
$ws.Range("D4").Value = 0.8928571428571429
$ws.Range("C5").Value = '(''Knn5'', (1.109078419726502, 0.7128458097177036))'
$ws.Range("D5").Value = 0.8928571428571429
$ws.Range("C6").Value = '(''LDA'', (0.9495932130840735, 0.6753846338630222))'
$ws.Range("D6").Value = 0.8928571428571429
$ws.Range("C7").Value = '(''Knn1'', (0.9475446232632139, 0.6472733703791613))'
$ws.Range("D7").Value = 0.8928571428571429
$ws.Range("C8").Value = '(''NaiveBayes'', (0.5704556562012737, 0.5294468484269127))'
$ws.Range("D8").Value = 0.8928571428571429
